# "Ahora los puntajes no descienden de 0"
# Mark the "Puntajes >= 0" backlog item (row 7) and the row-4 item as
# completed ("ok" in the Status column), and move the active selection
# to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status column (C) gets "ok" for row 4 ...
$ws.Range("C4").Value = "ok"

# ... and for row 7 ("Puntajes >= 0"), which previously had no status cell.
$ws.Range("C7").Value = "ok"

# Update the selected cell shown in the saved sheet view.
$ws.Range("C6").Select()
